$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 17.586354
$ws.Range("H2").Value = 52.759062
$ws.Range("I2").Value = 0.2178245326054132
$ws.Range("J2").Value = 0.2178245326054132
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.798756666666666
$ws.Range("N2").Value = 26.39627
$ws.Range("O2").Value = 0.3278753542842912
$ws.Range("P2").Value = 0.3278753542842912
$ws.Range("Q2").Value = 154.73804949986
$ws.Range("R2").Value = 1392.64244549874
$ws.Range("S2").Value = 0.07141929579980999
$ws.Range("T2").Value = 0.07141929579981

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 17.586354
$ws.Range("H3").Value = 52.759062
$ws.Range("I3").Value = 0.2178245326054132
$ws.Range("J3").Value = 0.2178245326054132
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.814629
$ws.Range("N3").Value = 5.443887
$ws.Range("O3").Value = 0.06762002278384967
$ws.Range("P3").Value = 0.06762002278384967
$ws.Range("Q3").Value = 31.912707972666
$ws.Range("R3").Value = 287.214371753994
$ws.Range("S3").Value = 0.01472929985765944
$ws.Range("T3").Value = 0.01472929985765945

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 17.586354
$ws.Range("H4").Value = 52.759062
$ws.Range("I4").Value = 0.2178245326054132
$ws.Range("J4").Value = 0.2178245326054132
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.22229
$ws.Range("N4").Value = 48.66687
$ws.Range("O4").Value = 0.6045046229318591
$ws.Range("P4").Value = 0.6045046229318591
$ws.Range("Q4").Value = 285.29093463066
$ws.Range("R4").Value = 2567.61841167594
$ws.Range("S4").Value = 0.1316759369479437
$ws.Range("T4").Value = 0.1316759369479438

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 52.27042766666667
$ws.Range("H5").Value = 156.811283
$ws.Range("I5").Value = 0.6474213742983183
$ws.Range("J5").Value = 0.6474213742983183
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.798756666666666
$ws.Range("N5").Value = 26.39627
$ws.Range("O5").Value = 0.3278753542842912
$ws.Range("P5").Value = 0.3278753542842912
$ws.Range("Q5").Value = 459.9147739016011
$ws.Range("R5").Value = 4139.23296511441
$ws.Range("S5").Value = 0.2122735124692838
$ws.Range("T5").Value = 0.2122735124692838

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 52.27042766666667
$ws.Range("H6").Value = 156.811283
$ws.Range("I6").Value = 0.6474213742983183
$ws.Range("J6").Value = 0.6474213742983183
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.814629
$ws.Range("N6").Value = 5.443887
$ws.Range("O6").Value = 0.06762002278384967
$ws.Range("P6").Value = 0.06762002278384967
$ws.Range("Q6").Value = 94.85143388633567
$ws.Range("R6").Value = 853.6629049770211
$ws.Range("S6").Value = 0.04377864808080354
$ws.Range("T6").Value = 0.04377864808080354

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 52.27042766666667
$ws.Range("H7").Value = 156.811283
$ws.Range("I7").Value = 0.6474213742983183
$ws.Range("J7").Value = 0.6474213742983183
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.22229
$ws.Range("N7").Value = 48.66687
$ws.Range("O7").Value = 0.6045046229318591
$ws.Range("P7").Value = 0.6045046229318591
$ws.Range("Q7").Value = 847.9460360326901
$ws.Range("R7").Value = 7631.51432429421
$ws.Range("S7").Value = 0.3913692137482309
$ws.Range("T7").Value = 0.3913692137482309

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 10.87955133333333
$ws.Range("H8").Value = 32.638654
$ws.Range("I8").Value = 0.1347540930962685
$ws.Range("J8").Value = 0.1347540930962685
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.798756666666666
$ws.Range("N8").Value = 26.39627
$ws.Range("O8").Value = 0.3278753542842912
$ws.Range("P8").Value = 0.3278753542842912
$ws.Range("Q8").Value = 95.72652482450889
$ws.Range("R8").Value = 861.53872342058
$ws.Range("S8").Value = 0.04418254601519738
$ws.Range("T8").Value = 0.04418254601519739

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 10.87955133333333
$ws.Range("H9").Value = 32.638654
$ws.Range("I9").Value = 0.1347540930962685
$ws.Range("J9").Value = 0.1347540930962685
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.814629
$ws.Range("N9").Value = 5.443887
$ws.Range("O9").Value = 0.06762002278384967
$ws.Range("P9").Value = 0.06762002278384967
$ws.Range("Q9").Value = 19.74234935645534
$ws.Range("R9").Value = 177.681144208098
$ws.Range("S9").Value = 0.009112074845386672
$ws.Range("T9").Value = 0.009112074845386673

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 10.87955133333333
$ws.Range("H10").Value = 32.638654
$ws.Range("I10").Value = 0.1347540930962685
$ws.Range("J10").Value = 0.1347540930962685
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.22229
$ws.Range("N10").Value = 48.66687
$ws.Range("O10").Value = 0.6045046229318591
$ws.Range("P10").Value = 0.6045046229318591
$ws.Range("Q10").Value = 176.49123679922
$ws.Range("R10").Value = 1588.42113119298
$ws.Range("S10").Value = 0.08145947223568441
$ws.Range("T10").Value = 0.08145947223568442

